{"js": "// Update the exploratory-questions bullet list to its final wording.\n//\n// Before (bulleted, numId=2):\n//   1. How much does a game's rating affect[ its sales?]   (2 runs)\n//   2. How much does a game's review (critical/user) affect its sales?\n//   3. Do different consoles sell more games than others?\n//   4. How do different regions affect game sales?\n//   5. Do certain genres sell more than others?\n//   6. How are genres trending over time?\n//\n// After:\n//   1. How much does a game's review (critical/user) affect its sales?\n//   2. Do different consoles sell more games than others?\n//   3. How do different regions affect game sales?\n//   4. Do certain genres sell more than others?\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the six consecutive bullet paragraphs by their current text so the\n// script is resilient to exact paragraph-index assumptions.\nconst items = paragraphs.items;\nlet startIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t.startsWith(\"How much does a game\") && t.indexOf(\"rating affect\") !== -1) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1) {\n  throw new Error(\"Could not locate the exploratory-questions bullet list.\");\n}\n\nconst targets = items.slice(startIndex, startIndex + 6);\n\nconst newText = [\n  \"How much does a game\\u2019s review (critical/user) affect its sales?\",\n  \"Do different consoles sell more games than others?\",\n  \"How do different regions affect game sales?\",\n  \"Do certain genres sell more than others?\",\n];\n\n// First four bullets get their text rewritten (this also merges the\n// first bullet's two runs back into a single run).\nfor (let i = 0; i < newText.length; i++) {\n  targets[i].insertText(newText[i], \"Replace\");\n}\n\n// Last two bullets (\"Do certain genres...\" / \"How are genres trending...\")\n// are removed outright.\ntargets[4].delete();\ntargets[5].delete();\n\nawait context.sync();\n", "ps1": "# Update the exploratory-questions bullet list to its final wording.\n#\n# Before (bulleted, numId=2):\n#   1. How much does a game's rating affect[ its sales?]   (2 runs)\n#   2. How much does a game's review (critical/user) affect its sales?\n#   3. Do different consoles sell more games than others?\n#   4. How do different regions affect game sales?\n#   5. Do certain genres sell more than others?\n#   6. How are genres trending over time?\n#\n# After:\n#   1. How much does a game's review (critical/user) affect its sales?\n#   2. Do different consoles sell more games than others?\n#   3. How do different regions affect game sales?\n#   4. Do certain genres sell more than others?\n\n$d = $word.ActiveDocument\n\n# Replace the text of the first four bullets in place. Using Find/Execute on\n# the whole paragraph's wording (rather than poking Range.Text directly)\n# correctly collapses the first bullet's two runs into one, matching the\n# target markup.\n$replacements = @(\n    @{ Find = \"How much does a game\u2019s rating affect its sales?\"; Replace = \"How much does a game\u2019s review (critical/user) affect its sales?\" },\n    @{ Find = \"How much does a game\u2019s review (critical/user) affect its sales?\"; Replace = \"Do different consoles sell more games than others?\" },\n    @{ Find = \"Do different consoles sell more games than others?\"; Replace = \"How do different regions affect game sales?\" },\n    @{ Find = \"How do different regions affect game sales?\"; Replace = \"Do certain genres sell more than others?\" }\n)\n\n# Paragraph 1 is the title, so the four bullets above are paragraphs 2-5.\n# Scoping each Find/Execute to its own paragraph's Range (instead of\n# searching the whole document) keeps each replacement targeted at the\n# right bullet, even once earlier edits make two bullets' wording coincide.\n$bulletIndex = 2\nforeach ($r in $replacements) {\n    $rng = $d.Paragraphs.Item($bulletIndex).Range\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n    $bulletIndex++\n}\n\n# The trailing two bullets (\"Do certain genres...\" and \"How are genres\n# trending...\") are removed outright. Delete from the bottom up so earlier\n# indices stay valid.\n$d.Paragraphs.Item(7).Range.Delete()\n$d.Paragraphs.Item(6).Range.Delete()\n"}
